$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.485.19"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "1.729.05"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'246.26"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "'0.4829"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "'0.2665"
$ws.Range("E8").Value = "  -1.38%  "
$ws.Range("D9").Value = "'0.06219"
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("D10").Value = "1.730.57"
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("D11").Value = "'0.07066"
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("D13").Value = "'4.597"
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("D14").Value = "'0.6104"
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("D15").Value = "'77.32"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "26.478.33"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "'0.000007199"
$ws.Range("E19").Value = "  +4.38%  "
$ws.Range("D20").Value = "'11.54"
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("D21").Value = "1.950.80"
$ws.Range("E21").Value = "  -0.86%  "
$ws.Range("D22").Value = "'4.506"
$ws.Range("E22").Value = "  -2.95%  "
$ws.Range("D23").Value = "'8.772"
$ws.Range("E23").Value = "  -0.92%  "
$ws.Range("D24").Value = "'5.249"
$ws.Range("E24").Value = "  -2.42%  "
$ws.Range("D25").Value = "'137.81"
$ws.Range("E25").Value = "  +1.22%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'1.777"
$ws.Range("E27").Value = "  -2.31%  "
$ws.Range("D28").Value = "'108.17"
$ws.Range("E28").Value = "  +1.20%  "
$ws.Range("D29").Value = "'1.403"
$ws.Range("E29").Value = "  -2.14%  "
$ws.Range("D30").Value = "'3.976"
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("D31").Value = "'0.07984"
$ws.Range("E31").Value = "  +1.19%  "
$ws.Range("D32").Value = "'3.690"
$ws.Range("E32").Value = "  -1.82%  "
$ws.Range("D33").Value = "'0.04565"
$ws.Range("E33").Value = "  -1.44%  "
$ws.Range("D34").Value = "'2.617"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").Value = "'1.004"
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("D36").Value = "'0.6332"
$ws.Range("E36").Value = "  -1.77%  "
$ws.Range("D37").Value = "'0.8910"
$ws.Range("E37").Value = "  -5.04%  "
$ws.Range("D38").Value = "'2.017"
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("D39").Value = "'2.391"
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("D42").Value = "'101.45"
$ws.Range("E42").Value = "  -10.63%  "
$ws.Range("D43").Value = "'5.475"
$ws.Range("E43").Value = "  -4.90%  "
$ws.Range("D44").Value = "'0.3886"
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("E45").Value = "  +3.45%  "
$ws.Range("E46").Value = "  -2.94%  "
$ws.Range("E47").Value = "  +0.91%  "
$ws.Range("D48").Value = "'7.888"
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("D50").Value = "'1.250"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("D51").Value = "'51.48"
$ws.Range("E51").Value = "  -0.52%  "
